$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of updated values scraped on Sun Mar 26 19:16:27 UTC 2023.
# Each entry: row number, new Price (column D), new Volume(1h) (column E).
# A price value is quoted so Excel/COM stores it as literal text (matching
# the workbook's existing inlineStr cells) instead of silently parsing it
# into a number.
$updates = @(
    @{ Row = 2; D = '27.881.52'; E = '  +1.56%  ' }
    @{ Row = 3; D = '1.770.20'; E = '  +2.02%  ' }
    @{ Row = 4; D = '1.002'; E = '  -0.42%  ' }
    @{ Row = 5; D = '328.17'; E = '  +1.94%  ' }
    @{ Row = 7; D = '0.4472'; E = '  -3.08%  ' }
    @{ Row = 8; D = '0.3565'; E = '  +1.28%  ' }
    @{ Row = 9; D = '0.07448'; E = '  +1.32%  ' }
    @{ Row = 10; D = '42.02'; E = '  +0.39%  ' }
    @{ Row = 11; D = '1.097'; E = '  +1.64%  ' }
    @{ Row = 12; D = $null; E = '  -0.40%  ' }
    @{ Row = 13; D = '20.97'; E = '  +2.63%  ' }
    @{ Row = 14; D = '6.026'; E = '  +1.97%  ' }
    @{ Row = 15; D = '7.245'; E = '  +2.87%  ' }
    @{ Row = 16; D = '1.771.28'; E = '  +1.88%  ' }
    @{ Row = 17; D = '93.36'; E = '  +2.60%  ' }
    @{ Row = 18; D = '0.00001061'; E = '  +0.88%  ' }
    @{ Row = 19; D = '0.06439'; E = '  +0.70%  ' }
    @{ Row = 20; D = $null; E = '  -0.35%  ' }
    @{ Row = 21; D = '17.11'; E = '  +2.90%  ' }
    @{ Row = 22; D = $null; E = '  +1.01%  ' }
    @{ Row = 23; D = '27.924.85'; E = '  +1.48%  ' }
    @{ Row = 24; D = '11.30'; E = '  +1.89%  ' }
    @{ Row = 25; D = '2.112'; E = '  +0.66%  ' }
    @{ Row = 26; D = '162.86'; E = '  +0.14%  ' }
    @{ Row = 27; D = '20.38'; E = '  +2.62%  ' }
    @{ Row = 28; D = '1.977.01'; E = '  +2.04%  ' }
    @{ Row = 29; D = $null; E = '  +6.28%  ' }
    @{ Row = 30; D = '125.06'; E = '  +0.42%  ' }
    @{ Row = 31; D = '1.106'; E = '  +6.08%  ' }
    @{ Row = 32; D = '0.09188'; E = '  +0.10%  ' }
    @{ Row = 33; D = '5.610'; E = '  +3.67%  ' }
    @{ Row = 35; D = '11.88'; E = '  +2.47%  ' }
    @{ Row = 36; D = '0.02293'; E = '  +1.31%  ' }
    @{ Row = 37; D = '0.06097'; E = '  +1.82%  ' }
    @{ Row = 38; D = '0.2102'; E = '  +1.83%  ' }
    @{ Row = 39; D = '0.6329'; E = '  +1.51%  ' }
    @{ Row = 40; D = '4.960'; E = '  +0.95%  ' }
    @{ Row = 41; D = '1.187'; E = '  +0.26%  ' }
    @{ Row = 42; D = '1.394'; E = '  +1.33%  ' }
    @{ Row = 43; D = '7.901'; E = '  +2.66%  ' }
    @{ Row = 44; D = '13.24'; E = '  +1.50%  ' }
    @{ Row = 45; D = '3.744'; E = $null }
    @{ Row = 46; D = '0.5907'; E = '  +1.86%  ' }
    @{ Row = 47; D = '122.30'; E = '  +0.58%  ' }
    @{ Row = 48; D = '1.956'; E = '  +2.03%  ' }
    @{ Row = 49; D = '0.06901'; E = '  +1.10%  ' }
    @{ Row = 50; D = '1.137'; E = '  +1.46%  ' }
    @{ Row = 51; D = '72.88'; E = '  +2.46%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)  # column D = Price
        $looksNumeric = $u.D -match '^-?\d+(\.\d+)?$'
        if ($looksNumeric) {
            # Force text storage so "1.002" etc. is not reinterpreted as a number.
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E  # column E = Volume(1h)
    }
}

Write-Host "Updated $($updates.Count) rows of crypto price/volume data."